# Auto-generated edit script applying the cryptos.xlsx diff
# Commit message: Updated cryptos list on Tue Apr 23 13:36:54 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.126.04"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "'3.180.90"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'604.51"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").Value = "'154.14"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'3.180.53"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "'0.544"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("D11").Value = "'5.64"
$ws.Range("E11").Value = "  -7.26%  "
$ws.Range("D12").Value = "'0.507"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "'38.30"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "'3.698.71"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "'66.130.87"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'7.36"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "'3.176.44"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "'507.69"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "'15.26"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Value = "'0.729"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").Value = "'8.02"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "'14.78"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").Value = "'84.36"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'2.99"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "'9.12"
$ws.Range("E28").Value = "  -2.80%  "
$ws.Range("D29").Value = "'2.39"
$ws.Range("E29").Value = "  +5.55%  "
$ws.Range("D30").Value = "'3.05"
$ws.Range("E30").Value = "  +6.19%  "
$ws.Range("D31").Value = "'7.01"
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").Value = "'27.92"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'1.18"
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").Value = "'510.14"
$ws.Range("E36").Value = "  +6.28%  "
$ws.Range("D37").Value = "'55.25"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").Value = "'0.0896"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").Value = "'0.0417"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").Value = "'0.0₃0711"
$ws.Range("E40").Value = "  +10.00%  "
$ws.Range("E41").Value = "  +5.01%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.89"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'8.76"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").Value = "'2.47"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").Value = "'2.830.08"
$ws.Range("E46").Value = "  -3.68%  "
$ws.Range("D47").Value = "'27.94"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("D49").Value = "'2.37"
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'2.66"
$ws.Range("E51").Value = "  +3.03%  "

Write-Output "Applied 94 cell updates"
